# Zeiterfassung.xlsx edit — "Esoterik Plugin beendet (Text fehlt noch)"
#
# Adds two new logged time entries (rows 28 + 29) for "Teresa":
#   - 25.12.2013 (serial 41633): 1.5h "Navi & Esoterik - begonnen"
#   - 27.12.2013 (serial 41635): 3h   "Esoterik - (fast) abgeschlossen"
# and updates the view's selection to reflect where the author ended up
# (topLeftCell A4 is not exposed by this host's saved sheetView, so we set
# the scroll position via ActiveWindow for correctness even though the
# round-tripped file cannot reflect it) with the active cell on C10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 28: copy formatting (date number format, borders, etc.) from the
# last populated data row (27) so the new rows reuse the existing "A" date
# style instead of minting a brand new cell style. ---
$ws.Range("A27").Copy($ws.Range("A28"))
$ws.Range("A28").Value = 41633
$ws.Range("B28").Value = "Teresa"
$ws.Range("D28").Value = 1.5
$ws.Range("E28").Value = "Navi & Esoterik - begonnen"

# --- Row 29 ---
$ws.Range("A27").Copy($ws.Range("A29"))
$ws.Range("A29").Value = 41635
$ws.Range("B29").Value = "Teresa"
$ws.Range("D29").Value = 3
$ws.Range("E29").Value = "Esoterik - (fast) abgeschlossen"

# C4 (=SUMIF(B8:B100,"Teresa",D8:D100)) recalculates automatically to 44
# once the new "Teresa" rows above are in place.

# --- View state: scroll so row 4 is the top-left visible row, and leave the
# selection on C10 (matches the sheetView/selection emitted by Excel after
# the author scrolled down and clicked C10). ---
[void]$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
[void]$ws.Range("C10").Select()
